# Update log 2019-7-17 21:49:13
# Appends two more logged time-blocks (afternoon 14:00-18:00 and evening
# 19:00-22:00) for 2019-07-17 to the bottom of the existing daily log table,
# widens the "expected goal" column to fit the new text, and leaves the
# selection on the last cell that was filled in (D9), mirroring what a user
# would do by typing straight into the worksheet after the last used row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 - afternoon block
$ws.Range("B8").Value = "下午14:00-18:00"
$ws.Range("C8").Value = "了解实体关系抽取"
$ws.Range("D8").Value = "浅入"

# Row 9 - evening block
$ws.Range("B9").Value = "晚上19:00-22:00"
$ws.Range("C9").Value = "初步学习发布部署服务器上的项目"
$ws.Range("D9").Value = "浅入"

# The new goal text is longer than what used to fit in column C, so widen it
# (target stored width ~31.5 chars; column widths are pixel-quantized, and
# 30.8 is the input that lands closest to that stored value).
$ws.Columns.Item(3).ColumnWidth = 30.8

# Leave the selection where the user finished typing.
[void]$ws.Range("D9").Select()
